# Append one new data row (row 71) to the ModCounts sheet, following the
# same pattern as the existing rows: Date, Game, ModCount.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = 70
$newRow = 71

# New row's Date column ("2026/01/20") looks like a date string; force it to
# be written as literal text (matching the existing column A cells, which are
# plain text) instead of letting it be auto-parsed into a date serial number.
$ws.Cells.Item($newRow, 1).NumberFormat = "@"
$ws.Cells.Item($newRow, 1).Value = "2026/01/20"
$ws.Cells.Item($newRow, 1).Style = "Normal"

$ws.Cells.Item($newRow, 2).Value = "逃离鸭科夫"
$ws.Cells.Item($newRow, 3).Value = 1149

# Match the formatting (centered alignment) used by the rest of the data rows.
$ws.Range("A$newRow`:C$newRow").HorizontalAlignment = $ws.Range("A$lastRow`:C$lastRow").HorizontalAlignment
$ws.Range("A$newRow`:C$newRow").VerticalAlignment = $ws.Range("A$lastRow`:C$lastRow").VerticalAlignment
